$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program")

# --- Fix up formatting first (so new cells get the plain "s4"-style look
#     already used by most of the data rows, instead of Excel's blank
#     default). We do this by copying format from an already-"plain"
#     formatted cell (A4, which uses that exact style) onto every cell
#     that will need it, before writing the new values. ---

$ws.Range("A4").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("D6").PasteSpecial(-4122)

$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("D7").PasteSpecial(-4122)

$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("D8").PasteSpecial(-4122)

$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("C9").PasteSpecial(-4122)

$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("B10").PasteSpecial(-4122)

$ws.Range("C11").PasteSpecial(-4122)

$ws.Range("D12").PasteSpecial(-4122)

# B6 is no longer part of the positive-flow data set - remove it entirely.
$ws.Cells.Item(6, 2).Clear()

# --- Now (re)write the test-data table, keeping only the positive flow
#     scenarios and adding the Edit* single-field scenarios. ---

# Row 2 - ValidInputAllOne
$ws.Cells.Item(2, 1).Value = "ValidInputAllOne"
$ws.Cells.Item(2, 2).Value = "poland"
$ws.Cells.Item(2, 3).Value = "for other modules"
$ws.Cells.Item(2, 4).Value = "active"

# Row 3 - ValidInputAllTwo
$ws.Cells.Item(3, 1).Value = "ValidInputAllTwo"
$ws.Cells.Item(3, 2).Value = "africa"
$ws.Cells.Item(3, 3).Value = "to edit and delete"
$ws.Cells.Item(3, 4).Value = "active"

# Row 4 - ValidInputThree
$ws.Cells.Item(4, 1).Value = "ValidInputThree"
$ws.Cells.Item(4, 2).Value = "australia"
$ws.Cells.Item(4, 3).Value = "multi delete"
$ws.Cells.Item(4, 4).Value = "active"

# Row 5 - ValidInputFour
$ws.Cells.Item(5, 1).Value = "ValidInputFour"
$ws.Cells.Item(5, 2).Value = "canada"
$ws.Cells.Item(5, 3).Value = "extra one"
$ws.Cells.Item(5, 4).Value = "active"

# Row 6 - BlankProgramName
$ws.Cells.Item(6, 1).Value = "BlankProgramName"
$ws.Cells.Item(6, 3).Value = "undergrad"
$ws.Cells.Item(6, 4).Value = "active"

# Row 7 - InvalidProgramName
$ws.Cells.Item(7, 1).Value = "InvalidProgramName"
$ws.Cells.Item(7, 2).Value = "a"
$ws.Cells.Item(7, 3).Value = "undergrad"
$ws.Cells.Item(7, 4).Value = "active"

# Row 8 - BlankDescriptionName
$ws.Cells.Item(8, 1).Value = "BlankDescriptionName"
$ws.Cells.Item(8, 2).Value = "randomName"
$ws.Cells.Item(8, 4).Value = "active"

# Row 9 - BlankStatus
$ws.Cells.Item(9, 1).Value = "BlankStatus"
$ws.Cells.Item(9, 2).Value = "randomName"
$ws.Cells.Item(9, 3).Value = "undergrad"

# Row 10 - EditProgramNameOnly
$ws.Cells.Item(10, 1).Value = "EditProgramNameOnly"
$ws.Cells.Item(10, 2).Value = "aaaaEdited"

# Row 11 - EditProgramDescOnly
$ws.Cells.Item(11, 1).Value = "EditProgramDescOnly"
$ws.Cells.Item(11, 3).Value = "knjjkd edited"

# Row 12 - EditProgramStatusOnly
$ws.Cells.Item(12, 1).Value = "EditProgramStatusOnly"
$ws.Cells.Item(12, 4).Value = "active"

# --- A11 / A12 pick up distinct "pasted-in" look (white fill, left
#     aligned, plain black Arial) - apply it to A11 then clone onto A12. ---

$a11 = $ws.Cells.Item(11, 1)
$a11.HorizontalAlignment = -4131
$a11.Font.Color = 0
$a11.Font.Name = "Arial"
$a11.Interior.Color = 16777215

$a11.Copy()
$ws.Cells.Item(12, 1).PasteSpecial(-4122)
